$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.582.52'
$ws.Range("E2").Value = '  +1.45%  '
$ws.Range("D3").Value = '1.908.12'
$ws.Range("E3").Value = '  +3.57%  '
$ws.Range("E4").Value = '  +0.97%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.15'
$ws.Range("E5").Value = '  +5.75%  '
$ws.Range("E6").Value = '  +2.36%  '
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.49'
$ws.Range("E8").Value = '  +2.82%  '
$ws.Range("E9").Value = '  +3.03%  '
$ws.Range("E10").Value = '  +1.98%  '
$ws.Range("E11").Value = '  +1.50%  '
$ws.Range("D12").Value = '2.185.74'
$ws.Range("E12").Value = '  +3.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '12.49'
$ws.Range("E13").Value = '  +9.34%  '
$ws.Range("D14").Value = '1.904.96'
$ws.Range("E14").Value = '  +3.34%  '
$ws.Range("E15").Value = '  +2.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.85'
$ws.Range("E16").Value = '  +3.75%  '
$ws.Range("D17").Value = '35.612.84'
$ws.Range("E17").Value = '  +1.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '71.99'
$ws.Range("E18").Value = '  +2.93%  '
$ws.Range("D19").Value = '0.0₃0811'
$ws.Range("E19").Value = '  +2.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '243.86'
$ws.Range("E20").Value = '  +1.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.46'
$ws.Range("E21").Value = '  +2.30%  '
$ws.Range("E22").Value = '  +3.46%  '
$ws.Range("E23").Value = '  +0.80%  '
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.94'
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.14'
$ws.Range("E26").Value = '  +23.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.52'
$ws.Range("E27").Value = '  +8.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.00'
$ws.Range("E28").Value = '  +2.92%  '
$ws.Range("E29").Value = '  +1.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.967'
$ws.Range("E30").Value = '  +28.06%  '
$ws.Range("E31").Value = '  +3.77%  '
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.15'
$ws.Range("E34").Value = '  +5.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.72'
$ws.Range("E35").Value = '  +5.72%  '
$ws.Range("E36").Value = '  +3.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.30'
$ws.Range("E37").Value = '  +4.15%  '
$ws.Range("E38").Value = '  +4.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0206'
$ws.Range("E39").Value = '  +3.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '91.68'
$ws.Range("E40").Value = '  +1.97%  '
$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.44'
$ws.Range("E41").Value = '  +50.38%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.357.94'
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.30'
$ws.Range("E43").Value = '  +4.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0592'
$ws.Range("E44").Value = '  +12.08%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("E45").Value = '  +3.52%  '
$ws.Range("B46").Value = 'Gas'
$ws.Range("C46").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.91'
$ws.Range("E46").Value = '  +16.84%  '
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.66'
$ws.Range("E49").Value = '  +4.84%  '
$ws.Range("D50").Value = '2.092.75'
$ws.Range("E50").Value = '  +3.26%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0689'
$ws.Range("E51").Value = '  +2.79%  '
